$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "29.123.44"
$ws.Range("E2").Value = "  -0.94%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.834.93"
$ws.Range("E3").Value = "  -0.94%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "0.9987"

# Row 5 - BNB
$ws.Range("D5").Value = "240.46"
$ws.Range("E5").Value = "  -1.89%  "

# Row 6 - XRP
$ws.Range("D6").Value = "0.6641"
$ws.Range("E6").Value = "  -4.00%  "

# Row 7 - USDC
$ws.Range("D7").Value = "0.9995"
$ws.Range("E7").Value = "  -0.12%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "0.2944"
$ws.Range("E8").Value = "  -3.69%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -4.00%  "

# Row 10 - Solana
$ws.Range("D10").Value = "22.72"
$ws.Range("E10").Value = "  -3.14%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.07695"
$ws.Range("E11").Value = "  -0.68%  "

# Row 12 - was Polkadot, now WrappedEther (rows 12/13 swapped)
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.818.15"
$ws.Range("E12").Value = "  -1.81%  "

# Row 13 - was WrappedEther, now Polkadot
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "5.015"
$ws.Range("E13").Value = "  -2.39%  "

# Row 14 - Polygon
$ws.Range("D14").Value = "0.6725"
$ws.Range("E14").Value = "  -2.87%  "

# Row 15 - Litecoin
$ws.Range("D15").Value = "86.00"
$ws.Range("E15").Value = "  -5.27%  "

# Row 16 - Uniswap
$ws.Range("D16").Value = "6.192"
$ws.Range("E16").Value = "  -1.78%  "

# Row 17 - ShibaInu
$ws.Range("D17").Value = "0.000008217"
$ws.Range("E17").Value = "  -0.63%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "28.728.86"
$ws.Range("E18").Value = "  -2.32%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "227.83"
$ws.Range("E19").Value = "  -3.61%  "

# Row 20 - Avalanche
$ws.Range("E20").Value = "  -1.57%  "

# Row 21 - Dai
$ws.Range("D21").Value = "0.9984"
$ws.Range("E21").Value = "  -0.22%  "

# Row 22 - Chainlink
$ws.Range("D22").Value = "7.230"
$ws.Range("E22").Value = "  -5.37%  "

# Row 23 - BinanceUSD
$ws.Range("D23").Value = "0.9999"
$ws.Range("E23").Value = "  -0.12%  "

# Row 24 - Monero
$ws.Range("D24").Value = "160.26"
$ws.Range("E24").Value = "  +0.28%  "

# Row 25 - Cosmos
$ws.Range("D25").Value = "8.682"
$ws.Range("E25").Value = "  -2.83%  "

# Row 26 - Stellar
$ws.Range("D26").Value = "0.1398"
$ws.Range("E26").Value = "  -5.37%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "18.00"
$ws.Range("E27").Value = "  -1.09%  "

# Row 28 - PancakeSwap
$ws.Range("D28").Value = "1.503"
$ws.Range("E28").Value = "  -1.48%  "

# Row 30 - InternetComputer(DFINITY)
$ws.Range("D30").Value = "4.070"
$ws.Range("E30").Value = "  -1.44%  "

# Row 31 - Toncoin
$ws.Range("D31").Value = "1.189"
$ws.Range("E31").Value = "  -1.21%  "

# Row 32 - Hedera
$ws.Range("D32").Value = "0.05345"
$ws.Range("E32").Value = "  +2.51%  "

# Row 33 - was ImmutableX, now LidoDAOToken (rows 33/34 swapped)
$ws.Range("B33").Value = "LidoDAOToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D33").Value = "1.854"
$ws.Range("E33").Value = "  -0.74%  "

# Row 34 - was LidoDAOToken, now ImmutableX
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "0.7483"
$ws.Range("E34").Value = "  -3.04%  "

# Row 35 - ARBITRUM
$ws.Range("D35").Value = "1.131"
$ws.Range("E35").Value = "  -1.08%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  -0.46%  "

# Row 37 - Maker
$ws.Range("D37").Value = "1.317.81"
$ws.Range("E37").Value = "  -0.67%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  -3.22%  "

# Row 39 - MXToken
$ws.Range("D39").Value = "2.714"
$ws.Range("E39").Value = "  -0.27%  "

# Row 40 - TrustWalletToken
$ws.Range("D40").Value = "0.9212"
$ws.Range("E40").Value = "  -1.85%  "

# Row 41 - FraxShare
$ws.Range("D41").Value = "5.960"
$ws.Range("E41").Value = "  +2.81%  "

# Row 42 - PaxDollar
$ws.Range("D42").Value = "0.9980"
$ws.Range("E42").Value = "  -0.23%  "

# Row 43 - Quant
$ws.Range("D43").Value = "103.31"
$ws.Range("E43").Value = "  -2.48%  "

# Row 44 - XinFinNetwork
$ws.Range("D44").Value = "0.08042"
$ws.Range("E44").Value = "  +16.23%  "

# Row 45 - BabyDogeCoin
$ws.Range("D45").Value = "0.00000000126"
$ws.Range("E45").Value = "  +2.28%  "

# Row 46 - Mantle
$ws.Range("D46").Value = "0.5165"
$ws.Range("E46").Value = "  -1.14%  "

# Row 47 - RocketPoolETH
$ws.Range("D47").Value = "1.948.03"
$ws.Range("E47").Value = "  -2.72%  "

# Row 48 - Aave
$ws.Range("D48").Value = "63.79"
$ws.Range("E48").Value = "  +1.43%  "

# Row 49 - RenderToken
$ws.Range("E49").Value = "  -1.63%  "

# Row 50 - EnergySwap
$ws.Range("D50").Value = "9.219"
$ws.Range("E50").Value = "  -4.81%  "

# Row 51 - Cronos
$ws.Range("D51").Value = "0.05925"
$ws.Range("E51").Value = "  -0.60%  "
